# Updates the cryptos price list (column D) and volume-change list (column E)
# for Mon May 15 07:51:32 UTC 2023, matching the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.818.45"
$ws.Range("D3").Value = "1.856.10"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  -1.82%  "
$ws.Range("D5").Value = "320.95"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("E6").Value = "  -1.75%  "
$ws.Range("D7").Value = "0.4315"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("D8").Value = "0.3797"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "0.07407"
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").Value = "0.8858"
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("D11").Value = "21.73"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "1.868.07"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "6.771"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").Value = "5.496"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").Value = "0.07105"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "88.47"
$ws.Range("E16").Value = "  +5.63%  "
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "0.000009043"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "1.014"
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("D20").Value = "15.54"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "27.857.48"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").Value = "5.281"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  -1.69%  "
$ws.Range("D24").Value = "2.090.02"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "2.031"
$ws.Range("E25").Value = "  +4.00%  "
$ws.Range("D26").Value = "156.62"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").Value = "18.63"
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").Value = "2.029"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("D29").Value = "5.418"
$ws.Range("E29").Value = "  +2.16%  "
$ws.Range("D30").Value = "121.43"
$ws.Range("E30").Value = "  +3.50%  "
$ws.Range("D31").Value = "0.08974"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").Value = "1.241"
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("D33").Value = "0.7777"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.580"
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("D36").Value = "1.148"
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("D38").Value = "0.05324"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01970"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").Value = "2.882"
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("D41").Value = "0.5204"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").Value = "7.011"
$ws.Range("E42").Value = "  +2.34%  "
$ws.Range("D43").Value = "0.1684"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").Value = "8.799"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").Value = "110.82"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("D46").Value = "10.77"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").Value = "0.4751"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("D48").Value = "1.713"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").Value = "0.06528"
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("E50").Value = "  -1.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.880"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.51%  "
